$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.079.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").Value = "'1.957.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'247.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "'0.4882"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.88%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.2967"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06819"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'19.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "'106.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.99%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.935.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "'0.07741"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.406"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.7103"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.23%  "

$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "'282.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.15%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'30.988.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007749"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'13.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.185.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.544"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'6.586"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.942"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.56%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'168.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.189"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.93%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.1057"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.440"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.773"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.06%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.513"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.94%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7668"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.09%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.165"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02048"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.132"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.19%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.411"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.34%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8821"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.20%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'73.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.82%  "

$ws.Range("D43").Value = "'109.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4442"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'1.018.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.94%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.461"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.64%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.377"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "

$ws.Range("D49").Value = "'0.1265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.71%  "

$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.2585"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.88%  "
